# Weekly refresh: insert 3 new rows of price data (for the newest reporting
# date) above the existing Comercializadora del Agro de Limari / Frutilla
# block, shifting the previously-recorded weeks down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 690:692, pushing rows 690-714 down to 693-717.
$ws.Range("A690:T692").Insert(-4121)

# Row 690 - Especial
$ws.Cells.Item(690, 1).Value = 2
$ws.Cells.Item(690, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(690, 3).Value = "Coquimbo"
$ws.Cells.Item(690, 4).Value = 45267
$ws.Cells.Item(690, 5).Value = 4
$ws.Cells.Item(690, 6).Value = "Fruta"
$ws.Cells.Item(690, 7).Value = 100101
$ws.Cells.Item(690, 8).Value = "Berries"
$ws.Cells.Item(690, 9).Value = 100112025
$ws.Cells.Item(690, 10).Value = "Frutilla"
$ws.Cells.Item(690, 11).Value = "Sin especificar"
$ws.Cells.Item(690, 12).Value = "Especial"
$ws.Cells.Item(690, 13).Value = 400
$ws.Cells.Item(690, 14).Value = 15000
$ws.Cells.Item(690, 15).Value = 16000
$ws.Cells.Item(690, 16).Value = 15500
$ws.Cells.Item(690, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(690, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(690, 19).Value = 2214
$ws.Cells.Item(690, 20).Value = 7

# Row 691 - Primera
$ws.Cells.Item(691, 1).Value = 2
$ws.Cells.Item(691, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(691, 3).Value = "Coquimbo"
$ws.Cells.Item(691, 4).Value = 45267
$ws.Cells.Item(691, 5).Value = 4
$ws.Cells.Item(691, 6).Value = "Fruta"
$ws.Cells.Item(691, 7).Value = 100101
$ws.Cells.Item(691, 8).Value = "Berries"
$ws.Cells.Item(691, 9).Value = 100112025
$ws.Cells.Item(691, 10).Value = "Frutilla"
$ws.Cells.Item(691, 11).Value = "Sin especificar"
$ws.Cells.Item(691, 12).Value = "Primera"
$ws.Cells.Item(691, 13).Value = 600
$ws.Cells.Item(691, 14).Value = 12000
$ws.Cells.Item(691, 15).Value = 13000
$ws.Cells.Item(691, 16).Value = 12500
$ws.Cells.Item(691, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(691, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(691, 19).Value = 1786
$ws.Cells.Item(691, 20).Value = 7

# Row 692 - Segunda
$ws.Cells.Item(692, 1).Value = 2
$ws.Cells.Item(692, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(692, 3).Value = "Coquimbo"
$ws.Cells.Item(692, 4).Value = 45267
$ws.Cells.Item(692, 5).Value = 4
$ws.Cells.Item(692, 6).Value = "Fruta"
$ws.Cells.Item(692, 7).Value = 100101
$ws.Cells.Item(692, 8).Value = "Berries"
$ws.Cells.Item(692, 9).Value = 100112025
$ws.Cells.Item(692, 10).Value = "Frutilla"
$ws.Cells.Item(692, 11).Value = "Sin especificar"
$ws.Cells.Item(692, 12).Value = "Segunda"
$ws.Cells.Item(692, 13).Value = 500
$ws.Cells.Item(692, 14).Value = 9000
$ws.Cells.Item(692, 15).Value = 10000
$ws.Cells.Item(692, 16).Value = 9500
$ws.Cells.Item(692, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(692, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(692, 19).Value = 1357
$ws.Cells.Item(692, 20).Value = 7
